$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 0.259369
$ws.Range("N2").Value = 0.778107
$ws.Range("O2").Value = 0.0514155333512404
$ws.Range("P2").Value = 0.0514155333512404
$ws.Range("Q2").Value = 0.09284900107633334
$ws.Range("R2").Value = 0.8356410096870001
$ws.Range("S2").Value = 0.0514155333512404
$ws.Range("T2").Value = 0.0514155333512404

# Row 3 updates
$ws.Range("O3").Value = 0.5982999525231611
$ws.Range("P3").Value = 0.5982999525231611
$ws.Range("S3").Value = 0.5982999525231611
$ws.Range("T3").Value = 0.5982999525231611

# Row 4 updates
$ws.Range("M4").Value = 1.767033
$ws.Range("N4").Value = 5.301099
$ws.Range("O4").Value = 0.3502845141255985
$ws.Range("P4").Value = 0.3502845141255985
$ws.Range("Q4").Value = 0.632563062351
$ws.Range("R4").Value = 5.693067561159
$ws.Range("S4").Value = 0.3502845141255985
$ws.Range("T4").Value = 0.3502845141255985
